$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B45").Value = "Lisboa"
$ws.Range("G45").Value = "extra_info: {`"date`": {`"comment`": `"entra para a China`"}, `"value`": {`"comment`": `"@wikidata:Q597`", `"original`": `"Lisboa, Arroios`"}}"
$ws.Range("B57").Value = "Lisboa"
$ws.Range("G57").Value = "extra_info: {`"value`": {`"comment`": `"@wikidata:Q597`", `"original`": `"Lisboa, Arroios`"}}"
$ws.Range("B61").Value = "Viena"
$ws.Range("G61").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("G108").Value = "extra_info: {`"value`": {`"comment`": `"[Arroios Rodrigues, Hist\u00f3ria Da Companhia de Jesus Na Assist\u00eancia de Portugal v4:167] @wikidata:Q597`"}}"
$ws.Range("B109").Value = "Palermo"
$ws.Range("F109").Value = "Palermo"
$ws.Range("G109").Value = "extra_info: {`"value`": {`"comment`": `"[Sic\u00edlia] @wikidata:Q2656`"}}"
$ws.Range("B116").Value = "Messina"
$ws.Range("F116").Value = "Messina"
$ws.Range("G116").Value = "extra_info: {`"value`": {`"comment`": `"[Sic\u00edlia] @wikidata:Q13666`"}}"
$ws.Range("B130").Value = "Viena"
$ws.Range("G130").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B135").Value = "Viena"
$ws.Range("G135").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("G156").Value = "extra_info: {`"value`": {`"comment`": `"[Arroios Rodrigues, Hist\u00f3ria Da Companhia de Jesus Na Assist\u00eancia de Portugal v4:167] @wikidata:Q597`"}}"
$ws.Range("B177").Value = "Viena"
$ws.Range("G177").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B244").Value = "Lisboa"
$ws.Range("G244").Value = "extra_info: {`"value`": {`"comment`": `"@wikidata:Q597`", `"original`": `"Lisboa, Arroios`"}}"
$ws.Range("B257").Value = "Brno"
$ws.Range("G257").Value = "extra_info: {`"value`": {`"comment`": `"[Rep\u00fablica Checa] @wikidata:Q14960`"}}"
$ws.Range("F259").Value = "Palermo"
$ws.Range("F260").Value = "Palermo"
$ws.Range("F266").Value = "Messina"
$ws.Range("B314").Value = "Paris"
$ws.Range("G314").Value = "Será que vale mesmo a pena chamar-lhe entrada? extra_info: {`"value`": {`"comment`": `"(Montmartre) @wikidata:Q186115`"}}"
$ws.Range("B316").Value = "Paris"
$ws.Range("G316").Value = "Será que vale mesmo a pena chamar-lhe entrada? extra_info: {`"value`": {`"comment`": `"(Montmartre) @wikidata:Q186115`"}}"
$ws.Range("B332").Value = "Brno"
$ws.Range("G332").Value = "extra_info: {`"date`": {`"comment`": `"ou 16840918`"}, `"value`": {`"comment`": `"[Rep\u00fablica Checa] @wikidata:Q14960`"}}"
$ws.Range("B334").Value = "Viena"
$ws.Range("G334").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B339").Value = "Viena"
$ws.Range("F339").Value = "Viena"
$ws.Range("G339").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B380").Value = "Palermo"
$ws.Range("G380").Value = "extra_info: {`"value`": {`"comment`": `"[Sic\u00edlia] @wikidata:Q2656`"}}"
$ws.Range("G382").Value = "`"https://roccadellacultura.it/opere-e-contenuti/la-spezieria-dei-gesuiti-di-novellara`" extra_info: {`"value`": {`"comment`": `"@wikidata:Q111218`"}}"
$ws.Range("B386").Value = "Palermo"
$ws.Range("F386").Value = "Palermo"
$ws.Range("G386").Value = "extra_info: {`"value`": {`"comment`": `"[Sic\u00edlia] @wikidata:Q2656`"}}"
$ws.Range("B397").Value = "Viena"
$ws.Range("F397").Value = "Viena"
$ws.Range("G397").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B451").Value = "Viena"
$ws.Range("G451").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B572").Value = "Viena"
$ws.Range("G572").Value = "extra_info: {`"date`": {`"comment`": `"j\u00e1 padre`"}, `"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B573").Value = "Viena"
$ws.Range("G573").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B575").Value = "Brno"
$ws.Range("G575").Value = "extra_info: {`"value`": {`"comment`": `"[Rep\u00fablica Checa] @wikidata:Q14960`"}}"
$ws.Range("B583").Value = "Lisboa"
$ws.Range("G583").Value = "extra_info: {`"value`": {`"comment`": `"@wikidata:@wikidata:Q15041623`", `"original`": `"Lisboa, Arroios`"}}"
$ws.Range("B637").Value = "Brno"
$ws.Range("G637").Value = "extra_info: {`"value`": {`"comment`": `"[Rep\u00fablica Checa] @wikidata:Q14960`"}}"
$ws.Range("B640").Value = "Viena"
$ws.Range("G640").Value = "destinado à China em 1785, viagem abortada extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B648").Value = "Palermo"
$ws.Range("G648").Value = "antes entrara na ordem de Malta em 16120518 diz Dehergne «La première entrée — avant « l'âge de raison » — est une « E. de dévotion », désirée par ses parents, qui consacrent leur fils au Seigneur et le revêtent des habits de l'Ordre, les Chevaliers de Malte.» extra_info: {`"value`": {`"comment`": `"[Sic\u00edlia] @wikidata:Q2656`"}}"
$ws.Range("F669").Value = "Caltavuturo, Palermo"
$ws.Range("B778").Value = "Viena"
$ws.Range("F778").Value = "Viena"
$ws.Range("G778").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("B797").Value = "Messina"
$ws.Range("G797").Value = "extra_info: {`"value`": {`"comment`": `"[Sic\u00edlia] @wikidata:Q13666`"}}"
$ws.Range("B855").Value = "Viena"
$ws.Range("G855").Value = "extra_info: {`"value`": {`"comment`": `"[\u00c1ustria] @wikidata:Q1741`"}}"
$ws.Range("E898").Value = "França"
$ws.Range("E899").Value = "França"
$ws.Range("B909").Value = "Lisboa"
$ws.Range("G909").Value = "extra_info: {`"value`": {`"comment`": `"@wikidata:Q597`", `"original`": `"Lisboa, Arroios`"}}"
$ws.Range("B917").Value = "Shiuchow"
$ws.Range("G917").Value = "extra_info: {`"date`": {`"comment`": `"fim do m\u00eas`"}, `"value`": {`"original`": `"Chao-tcheou fou`"}}"
$ws.Range("B918").Value = "Shiuchow"
$ws.Range("G918").Value = "extra_info: {`"date`": {`"comment`": `"fim do m\u00eas`"}, `"value`": {`"original`": `"Chao-tcheou fou`"}}"
$ws.Range("B954").Value = "Novellara"
$ws.Range("G954").Value = "extra_info: {`"value`": {`"comment`": `"[Assume-se Novellara, nascido a 40km] @wikidata:Q111218`", `"original`": `"Novellario`"}}"
$ws.Range("B961").Value = "Lisboa"
$ws.Range("G961").Value = "extra_info: {`"date`": {`"comment`": `"destinado \u00e0 China`"}, `"value`": {`"comment`": `"@wikidata:Q597`", `"original`": `"Lisboa, Arroios`"}}"
$ws.Range("F964").Value = "Palermo"
$ws.Range("B966").Value = "Brno"
$ws.Range("G966").Value = "extra_info: {`"value`": {`"comment`": `"[Rep\u00fablica Checa] @wikidata:Q14960`"}}"
$ws.Range("B967").Value = "Brno"
$ws.Range("G967").Value = "extra_info: {`"value`": {`"comment`": `"[Rep\u00fablica Checa] @wikidata:Q14960`"}}"
